# Insert a new weekly record at row 80 (pushing all subsequent rows down by
# one, so the former row 80..177 become rows 81..178, and a brand new row
# lands at 80).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("80:80").Insert()

# New row 80 values (columns that stay constant for every record in this
# sheet - A, B, C, E, F, G, H, I, N, Q, R - are re-populated with the same
# values used throughout the table; only D, J, K, L, M, O, P hold genuinely
# new data for this entry).
$ws.Cells.Item(80, 1).Value = 3
$ws.Cells.Item(80, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(80, 3).Value = "Coquimbo"
$ws.Cells.Item(80, 4).Value = 44413
$ws.Cells.Item(80, 5).Value = 5
$ws.Cells.Item(80, 6).Value = 100112031
$ws.Cells.Item(80, 7).Value = "Poroto verde"
$ws.Cells.Item(80, 8).Value = "Magnum"
$ws.Cells.Item(80, 9).Value = "Primera"
$ws.Cells.Item(80, 10).Value = 70
$ws.Cells.Item(80, 11).Value = 29000
$ws.Cells.Item(80, 12).Value = 30000
$ws.Cells.Item(80, 13).Value = 29500
$ws.Cells.Item(80, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(80, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(80, 16).Value = 1180
$ws.Cells.Item(80, 17).Value = 25
$ws.Cells.Item(80, 18).Value = "Hortaliza"
